$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.256.90"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.856.57"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "0.7026"
$ws.Range("E5").Value = "  +2.44%  "

$ws.Range("D6").Value = "237.93"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "0.08005"
$ws.Range("E8").Value = "  +4.40%  "

$ws.Range("D9").Value = "0.3025"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "23.58"
$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("D11").Value = "0.08192"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "5.196"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.7066"
$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.778.07"
$ws.Range("E14").Value = "  -5.21%  "

$ws.Range("D15").Value = "89.70"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "29.245.66"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "5.831"
$ws.Range("E17").Value = "  +1.87%  "

$ws.Range("D18").Value = "0.000007857"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").Value = "236.75"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.096.21"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "7.511"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").Value = "163.08"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").Value = "8.893"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").Value = "18.11"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").Value = "1.910"
$ws.Range("E29").Value = "  -2.48%  "

$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "1.474"
$ws.Range("E31").Value = "  -0.60%  "

$ws.Range("D32").Value = "4.345"
$ws.Range("E32").Value = "  -3.54%  "

$ws.Range("D33").Value = "4.025"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("D34").Value = "0.05171"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "1.166"
$ws.Range("E35").Value = "  -1.52%  "

$ws.Range("D36").Value = "0.7138"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").Value = "0.9931"
$ws.Range("E37").Value = "  -3.44%  "

$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("D41").Value = "0.9317"
$ws.Range("E41").Value = "  +1.85%  "

$ws.Range("D42").Value = "1.148.14"
$ws.Range("E42").Value = "  +4.09%  "

$ws.Range("D43").Value = "5.992"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").Value = "0.4259"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("D45").Value = "70.08"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").Value = "102.97"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").Value = "0.5289"

$ws.Range("D49").Value = "1.741"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").Value = "1.981.25"
$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").Value = "9.170"
$ws.Range("E51").Value = "  +0.22%  "
